# ==========================================================================
# Update "620101" indicator tables: add a second data series (DINEM - MIDES)
# next to the existing MIDES-MEF-OPP series, extend the year range with 2023
# and 2022, refresh the metadata sheet (observaciones + new actualizacion row).
# ==========================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Metadata")

# Helper: write a value as TEXT (shared string) even when it looks numeric
# (e.g. a year like "2023"), without leaving a stray number-format style
# behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Data sheet -----------------------------------------------------------

# Header row: rename the existing value column and add a second series column
$ws1.Range("B1").Value = "MIDES-MEF-OPP"
$ws1.Range("C1").Value = "DINEM - MIDES"

Set-TextValue $ws1.Range("A2") "2023"
$ws1.Range("B2").Value = 1.4

Set-TextValue $ws1.Range("A3") "2022"
$ws1.Range("B3").Value = 1.3

Set-TextValue $ws1.Range("A4") "2021"
$ws1.Range("B4").Value = 1.2

Set-TextValue $ws1.Range("A5") "2020"
$ws1.Range("B5").Value = 1.3

Set-TextValue $ws1.Range("A6") "2019"
$ws1.Range("B6").Value = 1.4

Set-TextValue $ws1.Range("A7") "2018"
$ws1.Range("B7").Value = 1.5
$ws1.Range("C7").Value = 1.6

Set-TextValue $ws1.Range("A8") "2017"
$ws1.Range("B8").Value = 1.5
$ws1.Range("C8").Value = 1.6

Set-TextValue $ws1.Range("A9") "2016"
$ws1.Range("B9").Value = 1.4
$ws1.Range("C9").Value = 1.6

Set-TextValue $ws1.Range("A10") "2015"
$ws1.Range("B10").Value = 1.3
$ws1.Range("C10").Value = 1.5

Set-TextValue $ws1.Range("A11") "2014"
$ws1.Range("B11").Value = 1.4
$ws1.Range("C11").Value = 1.5

Set-TextValue $ws1.Range("A12") "2013"
$ws1.Range("B12").Value = 1.4
$ws1.Range("C12").Value = 1.5

Set-TextValue $ws1.Range("A13") "2012"
$ws1.Range("B13").Value = 1.3
$ws1.Range("C13").Value = 1.5

Set-TextValue $ws1.Range("A14") "2011"
$ws1.Range("B14").Value = 1.4
$ws1.Range("C14").Value = 1.5

Set-TextValue $ws1.Range("A15") "2010"
$ws1.Range("B15").Value = 1.4
$ws1.Range("C15").Value = 1.5

Set-TextValue $ws1.Range("A16") "2009"
$ws1.Range("B16").Value = 1.5
$ws1.Range("C16").Value = 1.6

Set-TextValue $ws1.Range("A17") "2008"
$ws1.Range("B17").Value = 1.3
$ws1.Range("C17").Value = 1.4

Set-TextValue $ws1.Range("A18") "2007"
$ws1.Range("B18").Value = 1.3
$ws1.Range("C18").Value = 1.4

Set-TextValue $ws1.Range("A19") "2006"
$ws1.Range("B19").Value = 1.4
$ws1.Range("C19").Value = 1.5

Set-TextValue $ws1.Range("A20") "2005"
$ws1.Range("B20").Value = 1.2
$ws1.Range("C20").Value = 1.3

Set-TextValue $ws1.Range("A21") "2004"
$ws1.Range("B21").ClearContents()
$ws1.Range("C21").Value = 1.4

Set-TextValue $ws1.Range("A22") "2003"
$ws1.Range("B22").ClearContents()
$ws1.Range("C22").Value = 1.4

Set-TextValue $ws1.Range("A23") "2002"
$ws1.Range("B23").ClearContents()
$ws1.Range("C23").Value = 1.4

Set-TextValue $ws1.Range("A24") "2001"
$ws1.Range("B24").ClearContents()
$ws1.Range("C24").Value = 1.6

Set-TextValue $ws1.Range("A25") "2000"
$ws1.Range("B25").ClearContents()
$ws1.Range("C25").Value = 1.2

Set-TextValue $ws1.Range("A26") "1999"
$ws1.Range("B26").ClearContents()
$ws1.Range("C26").Value = 1.2

Set-TextValue $ws1.Range("A27") "1998"
$ws1.Range("B27").ClearContents()
$ws1.Range("C27").Value = 1.1

Set-TextValue $ws1.Range("A28") "1997"
$ws1.Range("B28").ClearContents()
$ws1.Range("C28").Value = 1.2

Set-TextValue $ws1.Range("A29") "1996"
$ws1.Range("B29").ClearContents()
$ws1.Range("C29").Value = 1.2

Set-TextValue $ws1.Range("A30") "1995"
$ws1.Range("B30").ClearContents()
$ws1.Range("C30").Value = 0.7

Set-TextValue $ws1.Range("A31") "1994"
$ws1.Range("B31").ClearContents()
$ws1.Range("C31").Value = 0.8

Set-TextValue $ws1.Range("A32") "1993"
$ws1.Range("B32").ClearContents()
$ws1.Range("C32").Value = 0.7

Set-TextValue $ws1.Range("A33") "1992"
$ws1.Range("B33").ClearContents()
$ws1.Range("C33").Value = 0.6

Set-TextValue $ws1.Range("A34") "1991"
$ws1.Range("B34").ClearContents()
$ws1.Range("C34").Value = 0.4

Set-TextValue $ws1.Range("A35") "1990"
$ws1.Range("B35").ClearContents()
$ws1.Range("C35").Value = 0.5

Set-TextValue $ws1.Range("A36") "1988"
$ws1.Range("B36").ClearContents()
$ws1.Range("C36").Value = 0.6

Set-TextValue $ws1.Range("A37") "1987"
$ws1.Range("B37").ClearContents()
$ws1.Range("C37").Value = 0.7

Set-TextValue $ws1.Range("A38") "1986"
$ws1.Range("B38").ClearContents()
$ws1.Range("C38").Value = 0.7

Set-TextValue $ws1.Range("A39") "1985"
$ws1.Range("B39").ClearContents()
$ws1.Range("C39").Value = 0.7

# --- Metadata sheet ---------------------------------------------------------

$ws2.Range("A1").Value = " "

# Expand the "observaciones" text to explain the two methodologies/series
$ws2.Range("B8").Value = "Las dos líneas representan metodologías ligeramente diferentes de cálculo. De acuerdo a lo informado en el Observatorio Social de MIDES, a partir del año 2016 se introdujo cambios en la metodología de estimación del Gasto Público Social producto de los cambios en la información brindada por el Presupuesto Nacional, lo cual llevó a trabajar en base al presupuesto por áreas programáticas (AP) de los incisos gubernamentales. El Gasto Público Social en Cultura y Deporte era considerado anteriormente bajo la denominación de Gasto Público Social No Convencional, definido como un subcomponente heterogéneo del GPS. La función Cultura y Deporte agrupa los gastos en museos, bibliotecas, organizaciones de prensa, servicios de televisión, deportes, y que antes también incluía otros conceptos que aludían a un aspecto multidisciplinario de los programas sociales. Se hizo una revisión de forma de dar consistencia en los conceptos para la serie desde 2015. La estimación siempre refiere a montos en pesos corrientes monto obligado intervenido por balance a partir de la información proporcionada mayoritariamente por Contaduría General de la Nación (CGN) del Ministerio de Economía y Finanzas (MEF). Para los años 2020 y 2021 se incluyen las erogaciones del fondo COVID destinadas a atender la emergencia sanitaria."

# Shift the trailing rows down to make room for the new "actualizacion" row
# (hardcode, since the COM layer's Range.Value getter is unreliable here —
# Value2 is used instead where a read-back is needed elsewhere)
$ws2.Range("A11").Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$ws2.Range("B11").Value = " "
$ws2.Range("A10").Value = "cita"
$ws2.Range("B10").Value = "UMAD con base en DINEM - MIDES hasta 2018, a partir de 2019 MIDES-MEF-OPP"

# New row: actualizacion / Julio 2025
$ws2.Range("A9").Value = "actualizacion"
$ws2.Range("B9").Value = "Julio 2025"

